$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" (column D) values for rows 2-49 (row -> new text).
$dValues = @{
    2  = "29.221.68"
    3  = "1.858.67"
    4  = "0.9996"
    5  = "0.7101"
    6  = "237.86"
    7  = "0.9997"
    8  = "0.07982"
    9  = "0.3035"
    10 = "23.59"
    11 = "0.08181"
    12 = "1.822.72"
    13 = "5.176"
    14 = "0.7040"
    15 = "89.73"
    16 = "29.236.72"
    17 = "5.847"
    18 = "0.000007883"
    19 = "13.27"
    20 = "238.02"
    21 = "0.9982"
    22 = "2.100.78"
    23 = "0.9996"
    24 = "7.423"
    25 = "162.39"
    26 = "8.933"
    27 = "0.1437"
    28 = "18.08"
    30 = "1.433"
    31 = "1.479"
    32 = "4.372"
    33 = "4.015"
    34 = "0.05212"
    35 = "1.162"
    36 = "0.7100"
    37 = "0.9982"
    38 = "2.669"
    39 = "0.01856"
    41 = "0.9300"
    42 = "1.128.90"
    43 = "0.4269"
    44 = "5.862"
    45 = "70.02"
    46 = "0.9987"
    47 = "102.98"
}

# New "Volume(1h)" (column E) values for rows 2-49 (row -> new text).
$eValues = @{
    2  = "  +0.34%  "
    3  = "  +0.31%  "
    4  = "  -0.12%  "
    5  = "  +2.17%  "
    6  = "  -0.36%  "
    7  = "  -0.09%  "
    8  = "  +4.47%  "
    9  = "  +0.07%  "
    10 = "  +0.92%  "
    11 = "  +0.67%  "
    12 = "  -1.75%  "
    13 = "  -1.18%  "
    14 = "  -3.10%  "
    15 = "  +0.68%  "
    16 = "  +0.40%  "
    17 = "  +1.07%  "
    18 = "  +1.64%  "
    19 = "  +0.71%  "
    20 = "  +0.54%  "
    21 = "  -0.19%  "
    22 = "  +0.13%  "
    23 = "  -0.16%  "
    24 = "  -2.49%  "
    25 = "  +0.59%  "
    26 = "  -0.70%  "
    27 = "  -0.61%  "
    28 = "  -0.04%  "
    29 = "  -2.26%  "
    30 = "  +1.97%  "
    31 = "  -1.00%  "
    32 = "  -2.39%  "
    33 = "  -0.14%  "
    34 = "  -0.33%  "
    35 = "  -2.30%  "
    36 = "  +1.19%  "
    37 = "  -0.94%  "
    38 = "  +0.58%  "
    39 = "  +0.04%  "
    40 = "  +1.83%  "
    41 = "  +0.01%  "
    42 = "  +4.38%  "
    43 = "  +0.02%  "
    44 = "  -2.78%  "
    45 = "  -0.60%  "
    46 = "  -0.20%  "
    47 = "  -0.09%  "
    48 = "  -4.14%  "
    49 = "  -0.59%  "
}

# Column D ("Price") holds values such as "29.246.65" or "237.92" that Excel
# would otherwise auto-convert to a number. Force the whole data range to
# Text before writing, then drop back to the default "Normal" style so no
# stray number-formatting is left behind on the cells.
$ws.Range("D2:D49").NumberFormat = "@"
foreach ($row in $dValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
}
$ws.Range("D2:D49").Style = "Normal"

foreach ($row in $eValues.Keys) {
    $ws.Cells.Item($row, 5).Value = $eValues[$row]
}

# Rows 50 and 51 swapped coins (RocketPoolETH now ranks above EnergySwap)
# along with their own updated price/volume figures.
$ws.Range("D50:D51").NumberFormat = "@"

$ws.Cells.Item(50, 2).Value = "RocketPoolETH"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(50, 4).Value = "1.985.71"
$ws.Cells.Item(50, 5).Value = "  -0.39%  "

$ws.Cells.Item(51, 2).Value = "EnergySwap"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(51, 4).Value = "9.158"
$ws.Cells.Item(51, 5).Value = "  -0.69%  "

$ws.Range("D50:D51").Style = "Normal"
